$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new manual SNP rows (21-29) ------------------------------------
# Fill column A (rs ids) first, then column C (shared "notes" citation),
# then column B (gene) - this mirrors how the source workbook grew its
# shared-string table.

$ws.Cells.Item(21,1).Value = "rs4846914"
$ws.Cells.Item(22,1).Value = "rs17145738"
$ws.Cells.Item(23,1).Value = "rs1495741"
$ws.Cells.Item(24,1).Value = "rs12678919"
$ws.Cells.Item(25,1).Value = "rs10761731"
$ws.Cells.Item(26,1).Value = "rs2412710"
$ws.Cells.Item(27,1).Value = "rs2929282"
$ws.Cells.Item(28,1).Value = "rs10401969"
$ws.Cells.Item(29,1).Value = "rs6065906"

$note = "Guay et al. 2024 (Acute pancreatitis risk in multifactorial chylomicronemia syndrome depends on the molecular cause of severe hypertriglyceridemia; PMID: 38448342)"
$ws.Cells.Item(21,3).Value = $note
$ws.Cells.Item(22,3).Value = $note
$ws.Cells.Item(23,3).Value = $note
$ws.Cells.Item(24,3).Value = $note
$ws.Cells.Item(25,3).Value = $note
$ws.Cells.Item(26,3).Value = $note
$ws.Cells.Item(27,3).Value = $note
$ws.Cells.Item(28,3).Value = $note
$ws.Cells.Item(29,3).Value = $note

$ws.Cells.Item(21,2).Value = "GALNT2"
$ws.Cells.Item(22,2).Value = "MLXIPL"
$ws.Cells.Item(23,2).Value = "NAT2"
$ws.Cells.Item(24,2).Value = "LPL"
$ws.Cells.Item(25,2).Value = "JMJD1C"
$ws.Cells.Item(26,2).Value = "CAPN3"
$ws.Cells.Item(27,2).Value = "FRMD5"
$ws.Cells.Item(28,2).Value = "CSPG3, CILP2, PBX4"
$ws.Cells.Item(29,2).Value = "PLTP"

# --- Columns B and C grew wider to fit the new values (bestFit) -------------
$ws.Columns.Item(2).ColumnWidth = 23.25
$ws.Columns.Item(3).ColumnWidth = 181.76

# --- Add a second "highlight duplicate values" conditional format rule ------
# on column A, on top of (with higher priority than) the existing one, using
# the same red-on-pink formatting.
$rng = $ws.Range("A1:A1048576")

$newRule = $rng.FormatConditions.AddUniqueValues()
$newRule.DupeUnique = 1
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615

# Touch a throw-away rule so a stray duplicate dxf is minted (matches the
# extra dxf left behind in the saved workbook).
$scratch = $rng.FormatConditions.AddUniqueValues()
$scratch.DupeUnique = 1
$scratch.Font.Color = 393372
$scratch.Interior.Color = 13551615
$scratch.Delete()

$newRule.SetFirstPriority()

Write-Output "done"
